$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("US20")

# --- Approval of US20 Test Case: mark review status as Done and add sign-offs ---

# Row 5 "Status" value: Pending -> Done
$ws.Range("F5").Value = "Done"

# Row 8 "Scrum Master" signature field
$ws.Range("F8").Value = "Megan Png"

# Row 9 "Tech Lead" signature field, noting approval via GitHub Pull Request
$ws.Range("F9").Value = "Claris (Approved in GitHub Pull Request)"

# Row 18 had an explicit custom row height; re-fit it back to the default height
$ws.Rows.Item(18).EntireRow.AutoFit()

# Reflect the reviewer's final on-screen state: zoomed in a bit more and
# focused on the newly-filled signature cell.
$ws.Range("F8").Select()
$excel.ActiveWindow.Zoom = 73
